$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Insert a new title row at the top of the "Total Wild Abundance" sheet,
# pushing the existing header + data rows down by one (matches the other
# two sheets, which already have a title row above their header row).
$ws1.Rows.Item(1).Insert() | Out-Null
$ws1.Range("A1").Value = "TOTAL WILD ABUNDANCE"

# Match the title formatting used on the other sheets (bold heading style)
# by copying the format from the analogous title cell on sheet 2.
$ws2.Range("A1").Copy() | Out-Null
$ws1.Range("A1").PasteSpecial(-4122) | Out-Null

# Restore the selections recorded in the saved workbook.
$ws2.Activate() | Out-Null
$ws2.Range("E16").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("D13").Select() | Out-Null
